$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "RET-26498"
$ws.Range("C2").Value = "Akter Telecom"
$ws.Range("D2").Value = "Md. Akter Hosen"
$ws.Range("F2").Value = 1797865555

# Row 3
$ws.Range("B3").Value = "RET-33094"
$ws.Range("C3").Value = "Shekh Electronics & Varieteis Store"
$ws.Range("D3").Value = "Emon Ali Shahin"
$ws.Range("F3").Value = 1712192700

# Row 4
$ws.Range("B4").Value = "RET-07894"
$ws.Range("C4").Value = "Bhai Bhai Store"
$ws.Range("D4").Value = "Md. Khokon Ahmed  "
$ws.Range("F4").Value = 1911861374

# Row 5
$ws.Range("B5").Value = "RET-33092"
$ws.Range("C5").Value = "Mondol Mobile Center"
$ws.Range("D5").Value = "Md. Azizul Mondol  "
$ws.Range("F5").Value = 1725821212

# Re-entering a numeric value into these cells clears their quotePrefix
# formatting; restore the original format (copied from an untouched
# sibling cell further down the same column) without disturbing the value.
$ws.Range("F8").Copy()
$ws.Range("F2:F5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New cell G5 with border (left/right thin) and no fill
$ws.Range("G5").Value = "With Contact Number"
$ws.Range("G5").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("G5").Borders.Item(7).Weight = 2      # xlThin
$ws.Range("G5").Borders.Item(10).LineStyle = 1  # xlEdgeRight
$ws.Range("G5").Borders.Item(10).Weight = 2     # xlThin

# Column C needs to widen to fit the new, longer retail-name text.
$ws.Columns.Item(3).ColumnWidth = 31.6

$ws.Range("J18").Select()
